# Add a "Weight" column (F) to the account info sheet.
# Rows 2-22 get a plain numeric weight; rows 23-30 get a "=1/8" formula
# (F24:F30 written as one range so Excel stores them as a shared formula
# group, matching how row 23 stands alone as a non-shared formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("F1").Value = "Weight"

# Rows 2-10: weight 0.1
$ws.Range("F2:F10").Value = 0.1

# Rows 11-12: weight 0.05
$ws.Range("F11:F12").Value = 0.05

# Rows 13-22: weight 0.1
$ws.Range("F13:F22").Value = 0.1

# Row 23: standalone formula
$ws.Range("F23").Formula = "=1/8"

# Rows 24-30: shared formula group
$ws.Range("F24:F30").Formula = "=1/8"

# Update the active selection to match the edited range
[void]$ws.Range("F23:F30").Select()
